$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was M) -> now B
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9527027027027027
$ws.Range("C2").Value = 0.986013986013986
$ws.Range("D2").Value = 0.9690721649484536
$ws.Range("E2").Value = 143

# Row 3 (was B) -> now M
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.975
$ws.Range("C3").Value = 0.9176470588235294
$ws.Range("D3").Value = 0.9454545454545454
$ws.Range("E3").Value = 85

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.9605263157894737
$ws.Range("C4").Value = 0.9605263157894737
$ws.Range("D4").Value = 0.9605263157894737
$ws.Range("E4").Value = 0.9605263157894737

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.9638513513513514
$ws.Range("C5").Value = 0.9518305224187578
$ws.Range("D5").Value = 0.9572633552014995
$ws.Range("E5").Value = 228

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9610152916073968
$ws.Range("C6").Value = 0.9605263157894737
$ws.Range("D6").Value = 0.9602673506634439
$ws.Range("E6").Value = 228
